# "Generate Report for Handoff"
#
# A new handoff round was generated for the file
# acaa613e-f43b-4227-878b-1390c9a6bb23.md (the 4th tracked source file,
# row 5 in each per-language table). This refreshes its
# "Latest Handoff Datetime" on each language sheet, and the Overview
# sheet's "Latest HO Xliff Generate Date" picks up the newest of those
# timestamps.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: Latest Handoff Datetime for acaa613e-...md (row 5, column H)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("H5").Value = "2016-08-30 08:13:20"

# de-de sheet: Latest Handoff Datetime for acaa613e-...md (row 5, column H)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("H5").Value = "2016-08-30 08:13:32"

# Overview sheet: Latest HO Xliff Generate Date for acaa613e-...md (row 5, column G)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G5").Value = "2016-08-30 08:13:32"
